$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.60024387145424
$ws.Range("C2").Value = 10.28170194026331
$ws.Range("E2").Value = 13.05809924185188
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.685313114342073
$ws.Range("K2").Value = 8.9805283210148
$ws.Range("L2").Value = 9.98922244188811
$ws.Range("M2").Value = 14.27773489005413
$ws.Range("O2").Value = 27.80131881722183
$ws.Range("B3").Value = 12.3948816509258
$ws.Range("C3").Value = 10.29171505908981
$ws.Range("E3").Value = 13.09078083702593
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.6870752044541
$ws.Range("K3").Value = 8.822164057191994
$ws.Range("L3").Value = 9.996868356097464
$ws.Range("M3").Value = 14.24857421657771
$ws.Range("O3").Value = 27.92158678243339
$ws.Range("B4").Value = 12.26916191965046
$ws.Range("C4").Value = 10.29848687171166
$ws.Range("E4").Value = 13.11271530629265
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.688214318123477
$ws.Range("K4").Value = 8.724638165517616
$ws.Range("L4").Value = 10.0028750634725
$ws.Range("M4").Value = 14.23257989911087
$ws.Range("O4").Value = 28.00064791980876
$ws.Range("B5").Value = 12.21808834220293
$ws.Range("C5").Value = 10.30140377948515
$ws.Range("E5").Value = 13.12212367241442
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.688692939757432
$ws.Range("K5").Value = 8.68487199552453
$ws.Range("L5").Value = 10.00565332468904
$ws.Range("M5").Value = 14.22654714869165
$ws.Range("O5").Value = 28.03417711447595
$ws.Range("B6").Value = 12.20961905896564
$ws.Range("C6").Value = 10.30189764738487
$ws.Range("E6").Value = 13.12371431061722
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.688773286969598
$ws.Range("K6").Value = 8.67826888494622
$ws.Range("L6").Value = 10.00613462771105
$ws.Range("M6").Value = 14.22557484174697
$ws.Range("O6").Value = 28.03982378756246
$ws.Range("B7").Value = 12.26847239752533
$ws.Range("C7").Value = 10.29852557242302
$ws.Range("E7").Value = 13.11284028796206
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 3.688220714522718
$ws.Range("K7").Value = 8.724101894817881
$ws.Range("L7").Value = 10.00291119329994
$ws.Range("M7").Value = 14.23249656940088
$ws.Range("O7").Value = 28.00109479814093
$ws.Range("B8").Value = 12.5293958615709
$ws.Range("C8").Value = 10.28502531093699
$ws.Range("E8").Value = 13.06898031618525
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.685908841402555
$ws.Range("K8").Value = 8.92601368139549
$ws.Range("L8").Value = 9.991586772959385
$ws.Range("M8").Value = 14.26728666499034
$ws.Range("O8").Value = 27.84170430340071
$ws.Range("B9").Value = 13.04124219856434
$ws.Range("C9").Value = 10.2634784088791
$ws.Range("E9").Value = 12.99778195704375
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.681826959770647
$ws.Range("K9").Value = 9.317554455286452
$ws.Range("L9").Value = 9.97976494115108
$ws.Range("M9").Value = 14.35045149279996
$ws.Range("O9").Value = 27.5705498344043
$ws.Range("B10").Value = 13.41375713955705
$ws.Range("C10").Value = 10.25062239854167
$ws.Range("E10").Value = 12.95448822732231
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.679100521830116
$ws.Range("K10").Value = 9.599804954326077
$ws.Range("L10").Value = 9.977374034709934
$ws.Range("M10").Value = 14.42035324783993
$ws.Range("O10").Value = 27.39660773853511
$ws.Range("B11").Value = 13.58171772201462
$ws.Range("C11").Value = 10.2454135833417
$ws.Range("E11").Value = 12.93674733311559
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.677918764712104
$ws.Range("K11").Value = 9.726489009168965
$ws.Range("L11").Value = 9.977644079385014
$ws.Range("M11").Value = 14.45399242961342
$ws.Range("O11").Value = 27.32297273411406
$ws.Range("B12").Value = 13.64504603398623
$ws.Range("C12").Value = 10.2435325858307
$ws.Range("E12").Value = 12.9303099793883
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.677479632455139
$ws.Range("K12").Value = 9.774172030132497
$ws.Range("L12").Value = 9.977940667783427
$ws.Range("M12").Value = 14.46698878647971
$ws.Range("O12").Value = 27.29587968305455
$ws.Range("B13").Value = 13.63142028239733
$ws.Range("C13").Value = 10.24393363195195
$ws.Range("E13").Value = 12.93168389609762
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.677573835672198
$ws.Range("K13").Value = 9.763916175964701
$ws.Range("L13").Value = 9.977868164532728
$ws.Range("M13").Value = 14.4641784272464
$ws.Range("O13").Value = 27.30167946062353
$ws.Range("B14").Value = 13.58693361611557
$ws.Range("C14").Value = 10.24525700216584
$ws.Range("E14").Value = 12.93621210292422
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.677882469450449
$ws.Range("K14").Value = 9.73041796069378
$ws.Range("L14").Value = 9.977664590698966
$ws.Range("M14").Value = 14.45505651800715
$ws.Range("O14").Value = 27.32072791430527
$ws.Range("B15").Value = 13.55964675857716
$ws.Range("C15").Value = 10.24607950207457
$ws.Range("E15").Value = 12.93902231463003
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.678072605830375
$ws.Range("K15").Value = 9.709860395814864
$ws.Range("L15").Value = 9.97756517505745
$ws.Range("M15").Value = 14.44950247048335
$ws.Range("O15").Value = 27.33249867960782
$ws.Range("B16").Value = 13.40274524548681
$ws.Range("C16").Value = 10.25097563280466
$ws.Range("E16").Value = 12.95568693266929
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.679178926392826
$ws.Range("K16").Value = 9.591487642553368
$ws.Range("L16").Value = 9.977383630992144
$ws.Range("M16").Value = 14.41819131480902
$ws.Range("O16").Value = 27.40153057352084
$ws.Range("B17").Value = 13.3060645128717
$ws.Range("C17").Value = 10.25414268386344
$ws.Range("E17").Value = 12.96641037920174
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.679872575852931
$ws.Range("K17").Value = 9.518399624297901
$ws.Range("L17").Value = 9.977619373503796
$ws.Range("M17").Value = 14.39944963396447
$ws.Range("O17").Value = 27.44528688531632
$ws.Range("B18").Value = 13.25031866969358
$ws.Range("C18").Value = 10.25602450469713
$ws.Range("E18").Value = 12.97276212075034
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.680277054619805
$ws.Range("K18").Value = 9.476202524987695
$ws.Range("L18").Value = 9.977882816312473
$ws.Range("M18").Value = 14.38884357489246
$ws.Range("O18").Value = 27.47097123096404
$ws.Range("B19").Value = 13.23142229261823
$ws.Range("C19").Value = 10.25667201395376
$ws.Range("E19").Value = 12.97494430164018
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.680414951761588
$ws.Range("K19").Value = 9.461889388388395
$ws.Range("L19").Value = 9.977993998243505
$ws.Range("M19").Value = 14.38528257244378
$ws.Range("O19").Value = 27.47975624765796
$ws.Range("B20").Value = 13.31637102493506
$ws.Range("C20").Value = 10.25379931716631
$ws.Range("E20").Value = 12.9652498185225
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.679798165720964
$ws.Range("K20").Value = 9.526196723858174
$ws.Range("L20").Value = 9.977581052644357
$ws.Range("M20").Value = 14.40142679262681
$ws.Range("O20").Value = 27.44057544855366
$ws.Range("B21").Value = 13.600008343263
$ws.Range("C21").Value = 10.24486581770237
$ws.Range("E21").Value = 12.93487444131067
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.677791589295259
$ws.Range("K21").Value = 9.740265387549687
$ws.Range("L21").Value = 9.977719118820371
$ws.Range("M21").Value = 14.45772889881516
$ws.Range("O21").Value = 27.31511144667179
$ws.Range("B22").Value = 13.78375488602037
$ws.Range("C22").Value = 10.23956017913881
$ws.Range("E22").Value = 12.91665862726831
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.676528964259798
$ws.Range("K22").Value = 9.878465919602709
$ws.Range("L22").Value = 9.978941576566859
$ws.Range("M22").Value = 14.49602566906984
$ws.Range("O22").Value = 27.23772416925081
$ws.Range("B23").Value = 13.68585364173855
$ws.Range("C23").Value = 10.24234329593468
$ws.Range("E23").Value = 12.92623110215149
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.677198400168729
$ws.Range("K23").Value = 9.804875431109702
$ws.Range("L23").Value = 9.978185846847834
$ws.Range("M23").Value = 14.47545100012316
$ws.Range("O23").Value = 27.27860493947256
$ws.Range("B24").Value = 13.31171195338367
$ws.Range("C24").Value = 10.2539543630806
$ws.Range("E24").Value = 12.96577392647736
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.679831788813311
$ws.Range("K24").Value = 9.522672205722573
$ws.Range("L24").Value = 9.977597979009575
$ws.Range("M24").Value = 14.40053239294588
$ws.Range("O24").Value = 27.44270384368652
$ws.Range("B25").Value = 12.90313808203292
$ws.Range("C25").Value = 10.26878295776963
$ws.Range("E25").Value = 13.01545871061994
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.682883155174986
$ws.Range("K25").Value = 9.212392171714821
$ws.Range("L25").Value = 9.981854746333314
$ws.Range("M25").Value = 14.26878421600718
$ws.Range("O25").Value = 27.63946811656485
